# Updated cryptos list on Fri Feb  2 12:54:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (Volume 1h) updates -----------------------------------------
# These values always contain non-numeric characters ("+", "%", spaces) so
# Excel keeps them as plain text automatically.
$volumeUpdates = @{
    "E2"  = "  +2.30%  "
    "E3"  = "  +1.96%  "
    "E4"  = "  -0.03%  "
    "E5"  = "  +1.03%  "
    "E6"  = "  +6.63%  "
    "E7"  = "  +2.12%  "
    "E8"  = "  -0.02%  "
    "E9"  = "  +5.12%  "
    "E10" = "  +9.21%  "
    "E11" = "  +1.02%  "
    "E12" = "  +3.15%  "
    "E13" = "  +13.84%  "
    "E14" = "  +3.45%  "
    "E15" = "  +1.83%  "
    "E16" = "  +2.25%  "
    "E17" = "  +3.40%  "
    "E18" = "  +2.11%  "
    "E19" = "  +7.96%  "
    "E20" = "  +3.97%  "
    "E22" = "  +2.48%  "
    "E23" = "  +0.50%  "
    "E24" = "  +12.99%  "
    "E25" = "  +0.28%  "
    "E26" = "  -0.06%  "
    "E27" = "  +4.24%  "
    "E28" = "  +3.42%  "
    "E30" = "  +0.58%  "
    "E31" = "  +0.52%  "
    "E32" = "  -0.08%  "
    "E33" = "  +2.78%  "
    "E34" = "  +0.16%  "
    "E35" = "  +4.00%  "
    "E36" = "  +2.30%  "
    "E37" = "  +0.84%  "
    "E38" = "  +4.56%  "
    "E39" = "  +2.32%  "
    "E40" = "  +4.39%  "
    "E41" = "  +1.62%  "
    "E42" = "  +1.40%  "
    "E43" = "  +4.69%  "
    "E44" = "  -1.96%  "
    "E45" = "  +7.21%  "
    "E46" = "  +5.25%  "
    "E47" = "  +0.49%  "
    "E48" = "  +8.29%  "
    "E49" = "  +1.97%  "
    "E50" = "  +3.71%  "
    "E51" = "  +2.18%  "
}

foreach ($key in $volumeUpdates.Keys) {
    $ws.Range($key).Value = $volumeUpdates[$key]
}

# --- Column D (Price) updates that are still non-numeric-looking text -----
# (multiple "." separators, e.g. "43.123.50") - Excel leaves these as text.
$priceTextUpdates = @{
    "D2"  = "43.123.50"
    "D3"  = "2.314.94"
    "D15" = "2.674.00"
    "D16" = "2.321.25"
    "D18" = "43.030.49"
    "D42" = "1.986.71"
    "D49" = "2.546.78"
}

foreach ($key in $priceTextUpdates.Keys) {
    $ws.Range($key).Value = $priceTextUpdates[$key]
}

# --- Column D (Price) updates that look like plain decimal numbers --------
# Excel's type-inference would silently store these as numbers. Force text
# (as the source sheet always stores prices as text) by switching the
# cell to a text format while writing, then restoring the original
# "Normal" style so no stray formatting/style is left behind.
$priceNumericLookingUpdates = @{
    "D5"  = "302.19"
    "D6"  = "101.48"
    "D9"  = "0.514"
    "D10" = "36.12"
    "D11" = "0.0794"
    "D13" = "18.02"
    "D14" = "6.89"
    "D19" = "12.65"
    "D22" = "67.82"
    "D23" = "236.16"
    "D24" = "2.20"
    "D26" = "1.00"
    "D27" = "24.74"
    "D29" = "34.75"
    "D30" = "168.63"
    "D31" = "9.20"
    "D32" = "0.999"
    "D33" = "5.03"
    "D34" = "4.73"
    "D35" = "17.39"
    "D37" = "0.0694"
    "D40" = "1.79"
    "D44" = "2.24"
    "D45" = "10.28"
    "D46" = "2.91"
    "D47" = "17.67"
    "D48" = "56.39"
    "D51" = "4.56"
}

foreach ($key in $priceNumericLookingUpdates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $priceNumericLookingUpdates[$key]
    $cell.Style = "Normal"
}

# --- Rows 33 / 34: RenderToken and Filecoin swapped places -----------------
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
